$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 of data (date 2020-05-14, index 4)
$ws.Range("A5").Value = 4
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 43965
$ws.Range("C5").Value = 24.8294
$ws.Range("D5").Value = 26.860399999999998
$ws.Range("E5").Value = 68.430000000000007
$ws.Range("F5").Value = 61.47
$ws.Range("G5").Value = 54.54
$ws.Range("H5").Value = 32.1
$ws.Range("I5").Value = 106.3

# Update the selection/view
$ws.Range("E4:H5").Select()
